$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 606, shifting existing rows 606:687 down to 607:688
$ws.Rows(606).Insert()

# Populate the newly inserted row 606 with the new data record
$ws.Cells.Item(606, 1).Value = 6
$ws.Cells.Item(606, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(606, 3).Value = "Metropolitana"
$ws.Cells.Item(606, 4).Value = 45077
$ws.Cells.Item(606, 5).Value = 13
$ws.Cells.Item(606, 6).Value = 100112039
$ws.Cells.Item(606, 7).Value = "Ciboulette"
$ws.Cells.Item(606, 8).Value = "Sin especificar"
$ws.Cells.Item(606, 9).Value = "Primera"
$ws.Cells.Item(606, 10).Value = 850
$ws.Cells.Item(606, 11).Value = 900
$ws.Cells.Item(606, 12).Value = 1000
$ws.Cells.Item(606, 13).Value = 955
$ws.Cells.Item(606, 14).Value = "$/docena de atados"
$ws.Cells.Item(606, 15).Value = "Región Metropolitana"
$ws.Cells.Item(606, 16).Value = 318
$ws.Cells.Item(606, 17).Value = 3
$ws.Cells.Item(606, 18).Value = "Hortaliza"
